$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Pomc -> Mc5r -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1885723333333333
$ws.Range("H2").Value = 0.565717
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01518366666666667
$ws.Range("N2").Value = 0.045551
$ws.Range("O2").Value = 0.02370341769240456
$ws.Range("P2").Value = 0.02370341769240456
$ws.Range("Q2").Value = 0.002863219451888889
$ws.Range("R2").Value = 0.025768975067
$ws.Range("S2").Value = 0.02370341769240456
$ws.Range("T2").Value = 0.02370341769240456

# Row 3 (ECs -> Pomc -> Mc5r -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1885723333333333
$ws.Range("H3").Value = 0.565717
$ws.Range("O3").Value = 0.8698427334878488
$ws.Range("P3").Value = 0.8698427334878488
$ws.Range("Q3").Value = 0.1050713726993333
$ws.Range("R3").Value = 0.945642354294
$ws.Range("S3").Value = 0.8698427334878488
$ws.Range("T3").Value = 0.8698427334878488

# Row 4 (ECs -> Pomc -> Mc5r -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1885723333333333
$ws.Range("H4").Value = 0.565717
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.068191
$ws.Range("N4").Value = 0.204573
$ws.Range("O4").Value = 0.1064538488197466
$ws.Range("P4").Value = 0.1064538488197466
$ws.Range("Q4").Value = 0.01285893598233333
$ws.Range("R4").Value = 0.115730423841
$ws.Range("S4").Value = 0.1064538488197466
$ws.Range("T4").Value = 0.1064538488197466
